# Commit: "Added logger and ThreadLocalDriver"
#
# Functional data change in the loginData sheet's test row: the expected
# result for the login test case in row 4 flips from "Invalid" to "Valid",
# and the active selection moves from the whole row (A4:C4) to just the
# cell that was edited (C4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loginData")

# C4 held the shared string "Invalid" (expectedResult column) - update it
# to "Valid" to match the rest of the positive-path test data.
$ws.Range("C4").Value = "Valid"

# Move/narrow the active selection to the cell that changed.
$ws.Range("C4").Select()
